$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Colab: For machine learning model development." ->
#           "Colab: For data analysis and machine learning model development."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    " For machine learning model development.", $true, $false, $false, $false, $false,
    $true, 1, $false, " For data analysis and machine learning model development.", 2
) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: append two new reference entries (title + hyperlink, each
# surrounded by the same blank-paragraph spacing used by the existing
# reference entries) after the Techmonitor reference at the end of the doc.
# ---------------------------------------------------------------------------

# Creates a brand-new empty paragraph right after the document's current
# last paragraph and returns it.
function New-TrailingParagraph($doc) {
    $cnt = $doc.Paragraphs.Count
    $lastPara = $doc.Paragraphs($cnt)
    $rng = $lastPara.Range
    $rng.Collapse(0)
    $rng.InsertParagraphAfter()
    $newCnt = $doc.Paragraphs.Count
    return $doc.Paragraphs($newCnt)
}

# Overwrites the (empty) paragraph's content with exact OOXML so the
# resulting markup matches the reference-list formatting used elsewhere in
# the document (jc=left paragraph mark + color/sz/szCs run formatting).
function Fill-ParagraphXML($para, $innerXml) {
    $rng = $para.Range
    $rng.Collapse(0)
    $fullXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">' + $innerXml + '</w:p>'
    $rng.InsertXML($fullXml)
}

$blankInner = '<w:pPr><w:jc w:val="left"/><w:rPr><w:color w:val="231f20"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r>'

# Blank paragraph
$p = New-TrailingParagraph($d)
Fill-ParagraphXML $p $blankInner

# Blank paragraph
$p = New-TrailingParagraph($d)
Fill-ParagraphXML $p $blankInner

# "Data analysis online with csv master:" title paragraph
$p = New-TrailingParagraph($d)
$innerXml = '<w:pPr><w:jc w:val="left"/><w:rPr><w:color w:val="231f20"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="231f20"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Data analysis online with csv master:</w:t></w:r>'
Fill-ParagraphXML $p $innerXml

# Hyperlink paragraph: https://www.csv-master.com/en
$p = New-TrailingParagraph($d)
$rng = $p.Range
$rng.Collapse(0)
$d.Hyperlinks.Add($rng, "https://www.csv-master.com/en", $null, $null, "https://www.csv-master.com/en") | Out-Null
$innerXml = '<w:pPr><w:jc w:val="left"/><w:rPr><w:color w:val="231f20"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:hyperlink r:id="rId10"><w:r><w:rPr><w:color w:val="1155cc"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">https://www.csv-master.com/en</w:t></w:r></w:hyperlink><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r>'
Fill-ParagraphXML $p $innerXml

# Blank paragraph
$p = New-TrailingParagraph($d)
Fill-ParagraphXML $p $blankInner

# "Visdium:" title paragraph
$p = New-TrailingParagraph($d)
$innerXml = '<w:pPr><w:jc w:val="left"/><w:rPr><w:color w:val="231f20"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="231f20"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Visdium:</w:t></w:r>'
Fill-ParagraphXML $p $innerXml

# Hyperlink paragraph: https://vizdium.com/
$p = New-TrailingParagraph($d)
$rng = $p.Range
$rng.Collapse(0)
$d.Hyperlinks.Add($rng, "https://vizdium.com/", $null, $null, "https://vizdium.com/") | Out-Null
$innerXml = '<w:pPr><w:jc w:val="left"/><w:rPr><w:color w:val="231f20"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:hyperlink r:id="rId11"><w:r><w:rPr><w:color w:val="1155cc"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">https://vizdium.com/</w:t></w:r></w:hyperlink><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r>'
Fill-ParagraphXML $p $innerXml

# Final trailing blank paragraph
$p = New-TrailingParagraph($d)
Fill-ParagraphXML $p $blankInner

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
